$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy the formatting from the
# existing header cell (H1, "IP") so the new headers match the bold /
# centered / bordered header style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2 for row 2 (plain, unstyled numbers like the
# rest of row 2).
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
